# Excel date support + time parse modification + tests
#
# 1. Add a new "DateSheet" worksheet at the end of the workbook.
# 2. Populate it with date/time values (built as DATE()+TIME() formulas,
#    then frozen to plain numeric values) plus a couple of helper columns.
# 3. Apply the assorted date/time number formats used by the fixture.
# 4. Make DateSheet the active/selected sheet (mirrors CustomSheet losing
#    its "tabSelected"/"topLeftCell" markers in the diff).

$wb = $excel.ActiveWorkbook

# --- New worksheet, inserted after the last existing sheet -----------------
# (CustomSheet was the previously-selected sheet; it loses its
# "tabSelected"/"topLeftCell" markers automatically once DateSheet becomes
# the active sheet below.)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "DateSheet"

# --- Row 1: header cells carrying only number formats, no values -----------
$ws.Range("A1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

# --- Row 2 ------------------------------------------------------------------
$ws.Range("A2").NumberFormat = "dd\-mm\-yy\ h:mm;@"
$ws.Range("A2").Formula = "=DATE(2017,3,14)+TIME(13,0,0)"
$ws.Range("A2").Value2 = $ws.Range("A2").Value2
$ws.Range("B2").Value = 1
$ws.Range("C2").NumberFormat = "dd\-mm\-yy\ h:mm;@"

# --- Row 3 ------------------------------------------------------------------
$ws.Range("A3").NumberFormat = "dd\-mm\-yy\ h:mm;@"
$ws.Range("A3").Formula = "=DATE(2017,3,15)+TIME(12,59,59)+0.995/86400"
$ws.Range("A3").Value2 = $ws.Range("A3").Value2
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = "date"
$ws.Range("E3").Value = "points"

# --- Row 4 ------------------------------------------------------------------
$ws.Range("A4").NumberFormat = "dd\-mm\-yy\ h:mm;@"
$ws.Range("A4").Formula = "=DATE(2019,5,19)+TIME(14,0,0)"
$ws.Range("A4").Value2 = $ws.Range("A4").Value2
$ws.Range("B4").Value = 3

$ws.Range("D4").NumberFormat = "[$-409]dd\-mm\-yy\ h:mm\ AM/PM;@"
$ws.Range("D4").Formula = "=DATE(2017,3,14)+TIME(13,0,0)"
$ws.Range("D4").Value2 = $ws.Range("D4").Value2
$ws.Range("E4").Value = 101

# --- Row 5 ------------------------------------------------------------------
$ws.Range("D5").NumberFormat = "[$-409]dd\-mm\-yy\ h:mm\ AM/PM;@"
$ws.Range("D5").Formula = "=DATE(2017,3,15)+TIME(12,59,59)+0.995/86400"
$ws.Range("D5").Value2 = $ws.Range("D5").Value2
$ws.Range("E5").Value = 102

# --- Row 6 ------------------------------------------------------------------
$ws.Range("D6").NumberFormat = "[$-409]dd\-mm\-yy\ h:mm\ AM/PM;@"
$ws.Range("D6").Formula = "=DATE(2019,5,19)+TIME(14,0,0)"
$ws.Range("D6").Value2 = $ws.Range("D6").Value2
$ws.Range("E6").Value = 103

# --- Column widths (closest values reachable through this host's
#     character-width quantization; OOXML targets are 13.44140625 /
#     25.21875 / 15.5546875 with column A & D "best fit") -------------------
$ws.Columns.Item(1).ColumnWidth = 12.59
$ws.Columns.Item(3).ColumnWidth = 24.25
$ws.Columns.Item(4).ColumnWidth = 14.59

# --- Selection / activation --------------------------------------------------
$ws.Activate()
$ws.Range("C13").Select()
